$d = $word.ActiveDocument

$null = $d.Content.Find.Execute("17×92=", $true, $false, $false, $false, $false, $true, 1, $false, "39×65=", 2)
$null = $d.Content.Find.Execute("79×58=", $true, $false, $false, $false, $false, $true, 1, $false, "58×71=", 2)
$null = $d.Content.Find.Execute("29×87=", $true, $false, $false, $false, $false, $true, 1, $false, "78×66=", 2)
$null = $d.Content.Find.Execute("20×33=", $true, $false, $false, $false, $false, $true, 1, $false, "71×48=", 2)
$null = $d.Content.Find.Execute("74×49=", $true, $false, $false, $false, $false, $true, 1, $false, "32×45=", 2)
$null = $d.Content.Find.Execute("89×16=", $true, $false, $false, $false, $false, $true, 1, $false, "58×48=", 2)
$null = $d.Content.Find.Execute("22×73=", $true, $false, $false, $false, $false, $true, 1, $false, "42×74=", 2)
$null = $d.Content.Find.Execute("85×27=", $true, $false, $false, $false, $false, $true, 1, $false, "40×15=", 2)
$null = $d.Content.Find.Execute("96×21=", $true, $false, $false, $false, $false, $true, 1, $false, "67×16=", 2)
$null = $d.Content.Find.Execute("10×53=", $true, $false, $false, $false, $false, $true, 1, $false, "96×67=", 2)
$null = $d.Content.Find.Execute("21×48=", $true, $false, $false, $false, $false, $true, 1, $false, "27×46=", 2)
$null = $d.Content.Find.Execute("86×98=", $true, $false, $false, $false, $false, $true, 1, $false, "43×40=", 2)
$null = $d.Content.Find.Execute("25×53=", $true, $false, $false, $false, $false, $true, 1, $false, "50×45=", 2)
$null = $d.Content.Find.Execute("66×93=", $true, $false, $false, $false, $false, $true, 1, $false, "63×97=", 2)
$null = $d.Content.Find.Execute("49×71=", $true, $false, $false, $false, $false, $true, 1, $false, "56×100=", 2)
$null = $d.Content.Find.Execute("30×18=", $true, $false, $false, $false, $false, $true, 1, $false, "44×100=", 2)
$null = $d.Content.Find.Execute("68×57=", $true, $false, $false, $false, $false, $true, 1, $false, "92×89=", 2)
$null = $d.Content.Find.Execute("14×68=", $true, $false, $false, $false, $false, $true, 1, $false, "70×50=", 2)
$null = $d.Content.Find.Execute("76×93=", $true, $false, $false, $false, $false, $true, 1, $false, "45×69=", 2)
$null = $d.Content.Find.Execute("61×45=", $true, $false, $false, $false, $false, $true, 1, $false, "25×76=", 2)
$null = $d.Content.Find.Execute("48×70=", $true, $false, $false, $false, $false, $true, 1, $false, "81×31=", 2)
$null = $d.Content.Find.Execute("46×34=", $true, $false, $false, $false, $false, $true, 1, $false, "55×10=", 2)
$null = $d.Content.Find.Execute("52×79=", $true, $false, $false, $false, $false, $true, 1, $false, "45×75=", 2)
$null = $d.Content.Find.Execute("10×17=", $true, $false, $false, $false, $false, $true, 1, $false, "72×96=", 2)
$null = $d.Content.Find.Execute("38×41=", $true, $false, $false, $false, $false, $true, 1, $false, "84×94=", 2)
$null = $d.Content.Find.Execute("56×25=", $true, $false, $false, $false, $false, $true, 1, $false, "41×95=", 2)
$null = $d.Content.Find.Execute("85×89=", $true, $false, $false, $false, $false, $true, 1, $false, "33×73=", 2)
$null = $d.Content.Find.Execute("64×41=", $true, $false, $false, $false, $false, $true, 1, $false, "64×23=", 2)
$null = $d.Content.Find.Execute("32×38=", $true, $false, $false, $false, $false, $true, 1, $false, "97×54=", 2)
$null = $d.Content.Find.Execute("42×21=", $true, $false, $false, $false, $false, $true, 1, $false, "38×93=", 2)
$null = $d.Content.Find.Execute("64×53=", $true, $false, $false, $false, $false, $true, 1, $false, "91×87=", 2)
$null = $d.Content.Find.Execute("26×72=", $true, $false, $false, $false, $false, $true, 1, $false, "95×11=", 2)
$null = $d.Content.Find.Execute("78×85=", $true, $false, $false, $false, $false, $true, 1, $false, "90×13=", 2)
$null = $d.Content.Find.Execute("97×67=", $true, $false, $false, $false, $false, $true, 1, $false, "96×42=", 2)
$null = $d.Content.Find.Execute("29×10=", $true, $false, $false, $false, $false, $true, 1, $false, "39×95=", 2)
$null = $d.Content.Find.Execute("96×89=", $true, $false, $false, $false, $false, $true, 1, $false, "35×23=", 2)
$null = $d.Content.Find.Execute("56×79=", $true, $false, $false, $false, $false, $true, 1, $false, "43×33=", 2)
$null = $d.Content.Find.Execute("80×66=", $true, $false, $false, $false, $false, $true, 1, $false, "34×16=", 2)
$null = $d.Content.Find.Execute("40×49=", $true, $false, $false, $false, $false, $true, 1, $false, "14×90=", 2)
$null = $d.Content.Find.Execute("95×46=", $true, $false, $false, $false, $false, $true, 1, $false, "50×51=", 2)
$null = $d.Content.Find.Execute("56×22=", $true, $false, $false, $false, $false, $true, 1, $false, "80×72=", 2)
$null = $d.Content.Find.Execute("57×21=", $true, $false, $false, $false, $false, $true, 1, $false, "10×88=", 2)
$null = $d.Content.Find.Execute("71×62=", $true, $false, $false, $false, $false, $true, 1, $false, "21×98=", 2)
$null = $d.Content.Find.Execute("42×89=", $true, $false, $false, $false, $false, $true, 1, $false, "42×95=", 2)
$null = $d.Content.Find.Execute("28×71=", $true, $false, $false, $false, $false, $true, 1, $false, "54×50=", 2)
$null = $d.Content.Find.Execute("39×15=", $true, $false, $false, $false, $false, $true, 1, $false, "70×94=", 2)
$null = $d.Content.Find.Execute("47×24=", $true, $false, $false, $false, $false, $true, 1, $false, "53×28=", 2)
$null = $d.Content.Find.Execute("28×92=", $true, $false, $false, $false, $false, $true, 1, $false, "23×94=", 2)
$null = $d.Content.Find.Execute("45×38=", $true, $false, $false, $false, $false, $true, 1, $false, "73×90=", 2)
$null = $d.Content.Find.Execute("68×63=", $true, $false, $false, $false, $false, $true, 1, $false, "45×87=", 2)
$null = $d.Content.Find.Execute("33×60=", $true, $false, $false, $false, $false, $true, 1, $false, "88×25=", 2)
$null = $d.Content.Find.Execute("72×67=", $true, $false, $false, $false, $false, $true, 1, $false, "14×93=", 2)
$null = $d.Content.Find.Execute("38×36=", $true, $false, $false, $false, $false, $true, 1, $false, "59×71=", 2)
$null = $d.Content.Find.Execute("44×22=", $true, $false, $false, $false, $false, $true, 1, $false, "37×42=", 2)
$null = $d.Content.Find.Execute("36×22=", $true, $false, $false, $false, $false, $true, 1, $false, "84×59=", 2)
$null = $d.Content.Find.Execute("56×74=", $true, $false, $false, $false, $false, $true, 1, $false, "45×97=", 2)
$null = $d.Content.Find.Execute("84×44=", $true, $false, $false, $false, $false, $true, 1, $false, "14×24=", 2)
$null = $d.Content.Find.Execute("90×12=", $true, $false, $false, $false, $false, $true, 1, $false, "40×29=", 2)
$null = $d.Content.Find.Execute("39×31=", $true, $false, $false, $false, $false, $true, 1, $false, "53×54=", 2)
$null = $d.Content.Find.Execute("74×26=", $true, $false, $false, $false, $false, $true, 1, $false, "54×28=", 2)
$null = $d.Content.Find.Execute("32×100=", $true, $false, $false, $false, $false, $true, 1, $false, "47×83=", 2)
$null = $d.Content.Find.Execute("13×19=", $true, $false, $false, $false, $false, $true, 1, $false, "69×67=", 2)
$null = $d.Content.Find.Execute("36×30=", $true, $false, $false, $false, $false, $true, 1, $false, "88×32=", 2)
$null = $d.Content.Find.Execute("46×43=", $true, $false, $false, $false, $false, $true, 1, $false, "62×87=", 2)
$null = $d.Content.Find.Execute("27×90=", $true, $false, $false, $false, $false, $true, 1, $false, "30×77=", 2)
$null = $d.Content.Find.Execute("49×100=", $true, $false, $false, $false, $false, $true, 1, $false, "66×95=", 2)
$null = $d.Content.Find.Execute("21×45=", $true, $false, $false, $false, $false, $true, 1, $false, "17×29=", 2)
$null = $d.Content.Find.Execute("22×57=", $true, $false, $false, $false, $false, $true, 1, $false, "84×84=", 2)
$null = $d.Content.Find.Execute("56×78=", $true, $false, $false, $false, $false, $true, 1, $false, "29×75=", 2)
$null = $d.Content.Find.Execute("37×38=", $true, $false, $false, $false, $false, $true, 1, $false, "59×99=", 2)
$null = $d.Content.Find.Execute("91×18=", $true, $false, $false, $false, $false, $true, 1, $false, "40×100=", 2)
$null = $d.Content.Find.Execute("20×58=", $true, $false, $false, $false, $false, $true, 1, $false, "52×64=", 2)
$null = $d.Content.Find.Execute("100×32=", $true, $false, $false, $false, $false, $true, 1, $false, "42×60=", 2)
$null = $d.Content.Find.Execute("36×25=", $true, $false, $false, $false, $false, $true, 1, $false, "91×17=", 2)
$null = $d.Content.Find.Execute("17×23=", $true, $false, $false, $false, $false, $true, 1, $false, "11×14=", 2)
$null = $d.Content.Find.Execute("20×64=", $true, $false, $false, $false, $false, $true, 1, $false, "61×67=", 2)
$null = $d.Content.Find.Execute("64×19=", $true, $false, $false, $false, $false, $true, 1, $false, "19×92=", 2)
$null = $d.Content.Find.Execute("67×66=", $true, $false, $false, $false, $false, $true, 1, $false, "57×48=", 2)
$null = $d.Content.Find.Execute("97×14=", $true, $false, $false, $false, $false, $true, 1, $false, "89×96=", 2)
$null = $d.Content.Find.Execute("55×18=", $true, $false, $false, $false, $false, $true, 1, $false, "32×75=", 2)
$null = $d.Content.Find.Execute("73×61=", $true, $false, $false, $false, $false, $true, 1, $false, "72×11=", 2)
$null = $d.Content.Find.Execute("81×89=", $true, $false, $false, $false, $false, $true, 1, $false, "76×57=", 2)
$null = $d.Content.Find.Execute("13×28=", $true, $false, $false, $false, $false, $true, 1, $false, "54×67=", 2)
$null = $d.Content.Find.Execute("50×97=", $true, $false, $false, $false, $false, $true, 1, $false, "42×49=", 2)
$null = $d.Content.Find.Execute("67×64=", $true, $false, $false, $false, $false, $true, 1, $false, "94×75=", 2)
$null = $d.Content.Find.Execute("91×42=", $true, $false, $false, $false, $false, $true, 1, $false, "43×43=", 2)
$null = $d.Content.Find.Execute("87×80=", $true, $false, $false, $false, $false, $true, 1, $false, "71×19=", 2)
$null = $d.Content.Find.Execute("72×15=", $true, $false, $false, $false, $false, $true, 1, $false, "31×85=", 2)
$null = $d.Content.Find.Execute("62×64=", $true, $false, $false, $false, $false, $true, 1, $false, "38×34=", 2)
$null = $d.Content.Find.Execute("27×91=", $true, $false, $false, $false, $false, $true, 1, $false, "98×16=", 2)
$null = $d.Content.Find.Execute("19×42=", $true, $false, $false, $false, $false, $true, 1, $false, "98×66=", 2)
$null = $d.Content.Find.Execute("86×36=", $true, $false, $false, $false, $false, $true, 1, $false, "81×90=", 2)
$null = $d.Content.Find.Execute("39×61=", $true, $false, $false, $false, $false, $true, 1, $false, "28×96=", 2)
$null = $d.Content.Find.Execute("11×99=", $true, $false, $false, $false, $false, $true, 1, $false, "63×52=", 2)
$null = $d.Content.Find.Execute("19×68=", $true, $false, $false, $false, $false, $true, 1, $false, "56×21=", 2)
$null = $d.Content.Find.Execute("83×60=", $true, $false, $false, $false, $false, $true, 1, $false, "79×92=", 2)
$null = $d.Content.Find.Execute("39×79=", $true, $false, $false, $false, $false, $true, 1, $false, "61×90=", 2)
$null = $d.Content.Find.Execute("83×48=", $true, $false, $false, $false, $false, $true, 1, $false, "69×51=", 2)
$null = $d.Content.Find.Execute("61×14=", $true, $false, $false, $false, $false, $true, 1, $false, "45×50=", 2)
$null = $d.Content.Find.Execute("65×16=", $true, $false, $false, $false, $false, $true, 1, $false, "77×95=", 2)
